# Auto-generated Excel COM-interop script
# Commit message: Updated cryptos list on Mon Oct 16 09:24:37 UTC 2023 with GitHub Actions
# Applies refreshed crypto market data (price/volume %) to the 'cryptos' sheet,
# including a rotation of three rows (26-28) and a swap of two rows (49-50)
# whose coin name/link moved to a different ranking position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to be treated as Text so that numeric-looking values
# (e.g. '210.80', '0.992') are NOT auto-converted to numbers by Excel, matching
# the workbook's original inline-string ('t="inlineStr"') cell type.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.455.41'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '1.566.05'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -1.29%  '
$ws.Range("D5").Value = '210.80'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").Value = '0.992'
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("D8").Value = '22.59'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").Value = '0.0869'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '1.789.24'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.562.16'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '27.443.95'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("D17").Value = '62.37'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '224.63'
$ws.Range("E18").Value = '  +3.85%  '
$ws.Range("D19").Value = '7.47'
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").Value = '0.992'
$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = '9.40'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("D25").Value = '149.82'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '6.60'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.107'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '15.10'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").Value = '0.993'
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").Value = '0.0470'
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '1.445.09'
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("D35").Value = '1.12'
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '0.539'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("D40").Value = '0.811'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").Value = '0.992'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = '1.84'
$ws.Range("E44").Value = '  +5.35%  '
$ws.Range("E45").Value = '  -2.98%  '
$ws.Range("D46").Value = '64.32'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = '1.701.43'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '86.58'
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0104'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0525'
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").Value = '0.0943'
$ws.Range("E51").Value = '  -1.91%  '

# Restore the default cell style now that the text values are safely stored,
# so no stray direct formatting is left behind on the edited cells.
$dataRange.Style = "Normal"
